$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing row 18..47 down to 19..48
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new data record
$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(18, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(18, 3).Value = "La Araucanía"
$ws.Cells.Item(18, 4).Value = 44495
$ws.Cells.Item(18, 5).Value = 9
$ws.Cells.Item(18, 6).Value = 300000001
$ws.Cells.Item(18, 7).Value = "Rabanito"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 30
$ws.Cells.Item(18, 11).Value = 6000
$ws.Cells.Item(18, 12).Value = 6000
$ws.Cells.Item(18, 13).Value = 6000
$ws.Cells.Item(18, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(18, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(18, 16).Value = 500
$ws.Cells.Item(18, 17).Value = 12
$ws.Cells.Item(18, 18).Value = "Hortaliza"
